# This script applies a weekly data update to the "Tuna" price sheet.
# A new week of price observations (date serial 45218) is inserted at the
# top of the data block (before the existing row 444), pushing all the
# existing data rows down by 3 rows. The three brand-new rows contain the
# latest week's quality/price data for "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 444; this shifts the existing
# rows 444:506 down to 447:509 and keeps all their original content/style.
$ws.Rows.Item(444).Resize(3).Insert()

# Data for the 3 new rows (444, 445, 446), same market/product metadata,
# new date + updated volume/price figures.
$newRows = @(
    @{ Row = 444; L = "Especial";               M = 280; N = 30000; O = 30000; P = 30000; S = 1667 },
    @{ Row = 445; L = "Extra (doble especial)";  M = 220; N = 32000; O = 32000; P = 32000; S = 1778 },
    @{ Row = 446; L = "Primera";                 M = 300; N = 25000; O = 25000; P = 25000; S = 1389 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 9
    $ws.Cells.Item($row, 2).Value  = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value  = "Metropolitana"
    $ws.Cells.Item($row, 4).Value  = 45218
    $ws.Cells.Item($row, 5).Value  = 13
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100107
    $ws.Cells.Item($row, 8).Value  = "Otros"
    $ws.Cells.Item($row, 9).Value  = 100107011
    $ws.Cells.Item($row, 10).Value = "Tuna"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/caja 18 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 18
}
